$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.726.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -6.86%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.299.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -8.16%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "185.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -10.71%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "519.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.61%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.598"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.298.47"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -8.02%  "

$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.626"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "60.74"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.00%  "

$ws.Range("E12").Value = "  -10.65%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000258"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.96%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.809.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.64%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.120"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.57%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.291.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -8.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -8.43%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "63.579.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.84%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.954"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -10.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "374.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -10.98%  "

$ws.Range("E26").Value = "  +1.50%  "

$ws.Range("E27").Value = "  -3.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.29%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.59%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "650.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -12.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.81"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.26"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "59.68"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.46%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.106"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.39%  "

$ws.Range("E37").Value = "  +0.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.394"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.54%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.66"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -11.85%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.997"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.21%  "

$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.974.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.96%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.126"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.48%  "

$ws.Range("E43").Value = "  -10.86%  "

$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -15.47%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.35%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0393"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.85%  "

$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.10%  "

$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.31%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.125"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -21.90%  "
